$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.278.59"
$ws.Range("E2").Value = "  -2.24%  "
$ws.Range("D3").Value = "1.559.76"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'206.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.20%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D8").Value = "'0.0612"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "'17.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "1.775.96"
$ws.Range("E12").Value = "  -3.65%  "
$ws.Range("D13").Value = "1.559.54"
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "25.278.31"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "'59.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").Value = "'1.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'186.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("D22").Value = "'9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'5.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "'141.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").Value = "'14.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").Value = "'6.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  -6.70%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("D33").Value = "'2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").Value = "'1.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "1.081.24"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("B38").Value = "PaxDollar"
$ws.Range("C38").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").Value = "'0.772"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.47%  "
$ws.Range("D42").Value = "'0.801"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.19%  "
$ws.Range("D43").Value = "'93.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Value = "1.691.00"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'52.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -2.07%  "
